# Q3 Update - 2025
# Applies quarterly data-refresh edits to the UN-MLS worksheet:
#  1) Updates the "short-url" value (column B) for every data row
#  2) Updates refreshed statistics (refugees/asylum_seekers/stateless/ooc)
#     for rows 963-1004 (columns N, O, S, T)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) short-url column (B2:B1005) changes from "5hWAV8" to "QDZ7yC" for all data rows
$ws.Range("B2:B1005").Value = "QDZ7yC"

# 2) Refreshed numeric-as-text statistics. The source column historically stores
#    these figures as text, so force Text number format before assigning so the
#    values are not reinterpreted as numbers.
$statCells = @("N963", "O963", "N964", "O964", "N965", "O965", "N966", "N967", "O967", "N969", "O969", "N970", "O970", "N972", "O973", "N974", "O974", "N975", "O975", "O977", "T978", "N979", "O979", "N980", "O980", "O982", "N983", "O983", "N985", "O985", "N987", "O987", "N988", "O989", "N990", "O990", "N991", "O991", "N993", "O993", "N994", "O994", "N995", "O995", "S996", "N997", "O997", "N998", "O998", "N999", "O999", "N1001", "O1001", "N1002", "O1003", "N1004", "O1004")
foreach ($addr in $statCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("N963").Value = "1193"
$ws.Range("O963").Value = "1555"
$ws.Range("N964").Value = "9"
$ws.Range("O964").Value = "0"
$ws.Range("N965").Value = "30"
$ws.Range("O965").Value = "50"
$ws.Range("N966").Value = "47"
$ws.Range("N967").Value = "74"
$ws.Range("O967").Value = "43"
$ws.Range("N969").Value = "16"
$ws.Range("O969").Value = "13"
$ws.Range("N970").Value = "9"
$ws.Range("O970").Value = "28"
$ws.Range("N972").Value = "35"
$ws.Range("O973").Value = "6"
$ws.Range("N974").Value = "194"
$ws.Range("O974").Value = "366"
$ws.Range("N975").Value = "10"
$ws.Range("O975").Value = "38"
$ws.Range("O977").Value = "10"
$ws.Range("T978").Value = "866"
$ws.Range("N979").Value = "238"
$ws.Range("O979").Value = "110"
$ws.Range("N980").Value = "200"
$ws.Range("O980").Value = "312"
$ws.Range("O982").Value = "6"
$ws.Range("N983").Value = "5"
$ws.Range("O983").Value = "5"
$ws.Range("N985").Value = "660"
$ws.Range("O985").Value = "445"
$ws.Range("N987").Value = "128365"
$ws.Range("O987").Value = "42169"
$ws.Range("N988").Value = "6"
$ws.Range("O989").Value = "50"
$ws.Range("N990").Value = "4135"
$ws.Range("O990").Value = "1574"
$ws.Range("N991").Value = "6"
$ws.Range("O991").Value = "17"
$ws.Range("N993").Value = "10"
$ws.Range("O993").Value = "27"
$ws.Range("N994").Value = "1999"
$ws.Range("O994").Value = "809"
$ws.Range("N995").Value = "64"
$ws.Range("O995").Value = "17"
$ws.Range("S996").Value = "120857"
$ws.Range("N997").Value = "227"
$ws.Range("O997").Value = "138"
$ws.Range("N998").Value = "401"
$ws.Range("O998").Value = "2281"
$ws.Range("N999").Value = "23"
$ws.Range("O999").Value = "13"
$ws.Range("N1001").Value = "14"
$ws.Range("O1001").Value = "7"
$ws.Range("N1002").Value = "5"
$ws.Range("O1003").Value = "5"
$ws.Range("N1004").Value = "359"
$ws.Range("O1004").Value = "2786"
